$wb = $excel.ActiveWorkbook

# The workbook contains one worksheet per "backward elimination" step; each
# sheet has the full OLS regression summary (statsmodels text output) stored
# in cell B2. That text embeds the run's Date: and Time: of generation which
# need to be refreshed on every sheet.
$oldDate = "Wed, 01 Jan 2020"
$newDate = "Thu, 02 Jan 2020"
$oldTime = "23:19:14"
$newTime = "20:49:08"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value()
    if ($text -ne $null -and ($text.Contains($oldDate) -or $text.Contains($oldTime))) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        $cell.Value = $updated
    }
}
